# Generate Report for Handoff
# Adds a new tracked file (d142da33-fb2d-4b61-88f4-a34802ac2cbf.md) as a new
# row at the bottom of each of the three tables (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newFile       = "d142da33-fb2d-4b61-88f4-a34802ac2cbf.md"
$newFilePath   = "e2e\d142da33-fb2d-4b61-88f4-a34802ac2cbf.md"
$overviewDate  = "2016-08-22 18:47:16"
$zhXlf         = "d142da33-fb2d-4b61-88f4-a34802ac2cbf.c202557b0636e130a850b77a13065b8c95c91b66.zh-cn.xlf"
$zhXlfDate     = "2016-08-22 18:47:09"
$deXlf         = "d142da33-fb2d-4b61-88f4-a34802ac2cbf.c202557b0636e130a850b77a13065b8c95c91b66.de-de.xlf"
$deXlfDate     = "2016-08-22 18:47:16"
$status        = "Ready for handoff"
$dateFormat    = "yyyy-mm-dd HH:mm:ss"
$srcRepoUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c202557b0636e130a850b77a13065b8c95c91b66/e2e/d142da33-fb2d-4b61-88f4-a34802ac2cbf.md"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOverview = $rowOverview.Range

$rngOverview.Cells.Item(1,1).Value = $newFile
$rngOverview.Cells.Item(1,2).Value = $newFilePath
$rngOverview.Cells.Item(1,3).Value = ".md"
$rngOverview.Cells.Item(1,4).Value = ""
$rngOverview.Cells.Item(1,5).Value = $status
$rngOverview.Cells.Item(1,6).Value = $status
$rngOverview.Cells.Item(1,7).Value = $overviewDate

$rngOverview.Cells.Item(1,7).NumberFormat = $dateFormat

$overviewRowIndex = $rowOverview.Range.Row
$bCell = $wsOverview.Cells.Item($overviewRowIndex, 2)
$wsOverview.Hyperlinks.Add($bCell, $srcRepoUrl, "", "", $newFilePath)

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$rngZh = $rowZh.Range

$rngZh.Cells.Item(1,1).Value = $newFile
$rngZh.Cells.Item(1,2).Value = ".md"
$rngZh.Cells.Item(1,3).Value = $status
$rngZh.Cells.Item(1,4).Value = "e2e"
$rngZh.Cells.Item(1,5).Value = "ht"
$rngZh.Cells.Item(1,6).Value = "False"
$rngZh.Cells.Item(1,7).Value = $zhXlf
$rngZh.Cells.Item(1,8).Value = $zhXlfDate
$rngZh.Cells.Item(1,9).Value = ""
$rngZh.Cells.Item(1,10).Value = ""
$rngZh.Cells.Item(1,11).Value = "0001-01-01 00:00:00"
$rngZh.Cells.Item(1,12).Value = ""
$rngZh.Cells.Item(1,13).Value = "True"
$rngZh.Cells.Item(1,14).Value = ""
$rngZh.Cells.Item(1,15).Value = "False"
$rngZh.Cells.Item(1,16).Value = ""

$rngZh.Cells.Item(1,8).NumberFormat = $dateFormat
$rngZh.Cells.Item(1,11).NumberFormat = $dateFormat

$zhRowIndex = $rowZh.Range.Row
$aCellZh = $wsZh.Cells.Item($zhRowIndex, 1)
$wsZh.Hyperlinks.Add($aCellZh, $srcRepoUrl, "", "", $newFile)

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$rngDe = $rowDe.Range

$rngDe.Cells.Item(1,1).Value = $newFile
$rngDe.Cells.Item(1,2).Value = ".md"
$rngDe.Cells.Item(1,3).Value = $status
$rngDe.Cells.Item(1,4).Value = "e2e"
$rngDe.Cells.Item(1,5).Value = "ht"
$rngDe.Cells.Item(1,6).Value = "False"
$rngDe.Cells.Item(1,7).Value = $deXlf
$rngDe.Cells.Item(1,8).Value = $deXlfDate
$rngDe.Cells.Item(1,9).Value = ""
$rngDe.Cells.Item(1,10).Value = ""
$rngDe.Cells.Item(1,11).Value = "0001-01-01 00:00:00"
$rngDe.Cells.Item(1,12).Value = ""
$rngDe.Cells.Item(1,13).Value = "True"
$rngDe.Cells.Item(1,14).Value = ""
$rngDe.Cells.Item(1,15).Value = "False"
$rngDe.Cells.Item(1,16).Value = ""

$rngDe.Cells.Item(1,8).NumberFormat = $dateFormat
$rngDe.Cells.Item(1,11).NumberFormat = $dateFormat

$deRowIndex = $rowDe.Range.Row
$aCellDe = $wsDe.Cells.Item($deRowIndex, 1)
$wsDe.Hyperlinks.Add($aCellDe, $srcRepoUrl, "", "", $newFile)
